# Actualización automática 2025-08-22 13:35:10
# Applies the updated commission/sales figures for LINDAO ZUÑIGA BRYAN JOSE
# across the three report sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("M28").Value = 143.97
$ws1.Range("D45").Value = 91.58
$ws1.Range("D46").Value = 88.53
$ws1.Range("M58").Value = 1190.95
$ws1.Range("D59").Value = "3 de 57"
$ws1.Range("M59").Value = "6 de 57"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F28").Value = 290.97
$ws2.Range("F45").Value = 91.58
$ws2.Range("F46").Value = 183.09
$ws2.Range("F58").Value = 1190.95
$ws2.Range("F59").Value = 82832.7

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D3").Value = 359.21
$ws3.Range("E3").Value = 5078.3732
$ws3.Range("F3").Value = 0.06606059839231516

$ws3.Range("D16").Value = 5747.39
$ws3.Range("E16").Value = 44642.78
$ws3.Range("F16").Value = 0.1140577616626418

$ws3.Range("D19").Value = 14353.84
$ws3.Range("E19").Value = 51001.26560036207
$ws3.Range("F19").Value = 0.219628441697759
